$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: add an explicit 8-decimal number format to the existing M17 value ---
$ws.Range("M17").NumberFormat = "0.00000000"

# --- Row 21 (new): weekly / daily row ---
$ws.Range("A21").Value = "weekly"
$ws.Range("B21").Value = "daily"
$ws.Range("C21").Value = 48.393621646
$ws.Range("C21").NumberFormat = "0.000000"
$ws.Range("D21").Value = 0.009157413
$ws.Range("E21").Value = 1.331405878
$ws.Range("E21").NumberFormat = "0.0000"

# --- Row 22 (new): weekly / weekly row ---
$ws.Range("A22").Value = "weekly"
$ws.Range("B22").Value = "weekly"
$ws.Range("C22").Value = 55.71193146
$ws.Range("C22").NumberFormat = "0.000000"
$ws.Range("D22").Value = 0.005775029
$ws.Range("D22").NumberFormat = "0.0000000"
$ws.Range("E22").Value = 1.134439089
$ws.Range("E22").NumberFormat = "0.0000"

# --- Row 23 (new): monthly / monthly row ---
$ws.Range("A23").Value = "monthly"
$ws.Range("B23").Value = "monthly"
$ws.Range("C23").Value = 64.844488952
$ws.Range("C23").NumberFormat = "0.000000"
$ws.Range("D23").Value = 0.002759109
$ws.Range("E23").Value = 0.845962132
$ws.Range("E23").NumberFormat = "0.0000"

# --- Row 25 (new, one blank row gap after row 23): summary values ---
# C25 reuses the look of the other "plain Calibri, vertically centered" value
# cells already in the sheet (e.g. K6) before getting its own number format.
$ws.Range("K6").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C25").Value = 24.26512
$ws.Range("C25").NumberFormat = "0.0"
$ws.Range("D25").Value = 0.0009255197
$ws.Range("D25").NumberFormat = "0.000000"
$ws.Range("E25").Value = 1.611847
$ws.Range("E25").NumberFormat = "0.0000"

# --- Selection matches the author's final cursor position ---
[void]$ws.Range("E25").Select()
